# edit.ps1 - applies the OOXML diff to before.docx via Word COM-interop
$d = $word.ActiveDocument

function Set-ParagraphXml {
    param(
        [int]$ParaIndex,
        [string]$InnerXml
    )
    $p = $d.Paragraphs($ParaIndex)
    $r = $p.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $InnerXml + '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

$rPrCommon = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="es-GT"/></w:rPr>'

# Paragraph 2: "Entidad XXXXXXX"
$inner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' + $rPrCommon + '</w:pPr>' +
         '<w:r>' + $rPrCommon + '<w:t>Entidad XXXXXXX</w:t></w:r></w:p>'
Set-ParagraphXml 2 $inner

# Paragraph 3: "Auditoria de Procesos" (accent removed)
$inner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' + $rPrCommon + '</w:pPr>' +
         '<w:r>' + $rPrCommon + '<w:t>Auditoria de Procesos</w:t></w:r></w:p>'
Set-ParagraphXml 3 $inner

# Paragraph 4: "Del 01 de Enero al 31 de Diciembre de 2024" (merged into one run)
$inner = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="center"/>' + $rPrCommon + '</w:pPr>' +
         '<w:r>' + $rPrCommon + '<w:t>Del 01 de Enero al 31 de Diciembre de 2024</w:t></w:r></w:p>'
Set-ParagraphXml 4 $inner

Write-Output "done"
